$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for price cells whose values look like plain numbers,
# matching the source file where these are inline strings (not numeric cells).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '42.165.53'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '2.225.09'
$ws.Range('E3').Value = '  -0.86%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '243.30'
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range('D6').Value = '0.628'
$ws.Range('E6').Value = '  +1.49%  '
$ws.Range('D7').Value = '73.25'
$ws.Range('E7').Value = '  -0.71%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').Value = '0.609'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').Value = '42.82'
$ws.Range('E10').Value = '  +3.54%  '
$ws.Range('D11').Value = '0.0964'
$ws.Range('E11').Value = '  +3.22%  '
$ws.Range('D12').Value = '7.07'
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('D13').Value = '0.104'
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').Value = '14.26'
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('D15').Value = '0.844'
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('D16').Value = '2.226.55'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = '42.046.97'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('D18').Value = '0.0000112'
$ws.Range('E18').Value = '  +14.98%  '
$ws.Range('D19').Value = '6.17'
$ws.Range('E19').Value = '  +1.77%  '
$ws.Range('D20').Value = '72.40'
$ws.Range('E20').Value = '  +1.06%  '
$ws.Range('D21').Value = '10.21'
$ws.Range('E21').Value = '  +38.70%  '
$ws.Range('D22').Value = '230.25'
$ws.Range('E22').Value = '  +0.49%  '
$ws.Range('D23').Value = '2.13'
$ws.Range('E23').Value = '  -9.34%  '
$ws.Range('D24').Value = '11.78'
$ws.Range('E24').Value = '  +6.71%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').Value = '3.59'
$ws.Range('E26').Value = '  +1.58%  '
$ws.Range('D27').Value = '2.29'
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('E28').Value = '  -0.64%  '
$ws.Range('D29').Value = '166.60'
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('D30').Value = '20.68'
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.0798'
$ws.Range('E31').Value = '  -3.87%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '5.53'
$ws.Range('E32').Value = '  +14.73%  '
$ws.Range('D33').Value = '0.116'
$ws.Range('E33').Value = '  -2.81%  '
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').Value = '0.124'
$ws.Range('E34').Value = '  -0.17%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').Value = '29.41'
$ws.Range('E35').Value = '  -2.44%  '
$ws.Range('D36').Value = '4.36'
$ws.Range('E36').Value = '  -2.62%  '
$ws.Range('D37').Value = '0.0302'
$ws.Range('E37').Value = '  +1.53%  '
$ws.Range('D38').Value = '13.11'
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('D39').Value = '2.15'
$ws.Range('E39').Value = '  -0.72%  '
$ws.Range('B40').Value = 'MultiversX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D40').Value = '64.88'
$ws.Range('E40').Value = '  +6.40%  '
$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').Value = '5.58'
$ws.Range('E41').Value = '  -3.27%  '
$ws.Range('D42').Value = '0.199'
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('D43').Value = '8.73'
$ws.Range('E43').Value = '  +0.74%  '
$ws.Range('D44').Value = '104.85'
$ws.Range('E44').Value = '  -4.64%  '
$ws.Range('D45').Value = '0.101'
$ws.Range('E45').Value = '  +1.41%  '
$ws.Range('D46').Value = '2.39'
$ws.Range('E46').Value = '  +6.30%  '
$ws.Range('D47').Value = '1.12'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').Value = '1.16'
$ws.Range('E48').Value = '  +0.75%  '
$ws.Range('E49').Value = '  +0.93%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.432.73'
$ws.Range('E50').Value = '  -0.77%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = '4.03'
$ws.Range('E51').Value = '  -1.31%  '
